$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.914.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.85%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.414.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.65%  '

$ws.Range("E10").Value = '  -1.46%  '

$ws.Range("E11").Value = '  +1.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.000.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.40%  '

$ws.Range("E13").Value = '  +0.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.430.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.23%  '

$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.940.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.565'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.556.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.45%  '

$ws.Range("E26").Value = '  -3.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("E33").Value = '  -2.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.82%  '

$ws.Range("E36").Value = '  +1.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '169.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.86'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '30.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.449.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("E41").Value = '  +1.30%  '

$ws.Range("E42").Value = '  +0.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.774'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.83%  '

$ws.Range("E44").Value = '  -1.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.534.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.53%  '

$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.00%  '
